$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.350.61"
$ws.Range("E2").Value = "  -2.53%  "
$ws.Range("D3").Value = "1.573.75"
$ws.Range("E3").Value = "  -3.60%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.17"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -2.92%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  -4.81%  "
$ws.Range("E8").Value = "  -2.10%  "
$ws.Range("E9").Value = "  -1.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.01"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -1.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0782"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").Value = "1.792.25"
$ws.Range("E12").Value = "  -3.56%  "
$ws.Range("D13").Value = "1.586.08"
$ws.Range("E13").Value = "  -2.81%  "
$ws.Range("E14").Value = "  -3.20%  "
$ws.Range("E15").Value = "  -3.45%  "
$ws.Range("D16").Value = "25.331.53"
$ws.Range("E16").Value = "  -2.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "59.84"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -2.51%  "
$ws.Range("D18").Value = "0.0₃0712"
$ws.Range("E18").Value = "  -3.83%  "
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "186.02"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -2.30%  "
$ws.Range("E21").Value = "  -2.10%  "
$ws.Range("E22").Value = "  -3.01%  "
$ws.Range("E23").Value = "  -2.52%  "
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.27"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -2.05%  "
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("E27").Value = "  -5.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "14.90"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.47"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -4.25%  "
$ws.Range("E30").Value = "  -6.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0461"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -3.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.06"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -2.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.02"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -2.99%  "
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.25"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -6.54%  "
$ws.Range("D36").Value = "1.086.87"
$ws.Range("E36").Value = "  -3.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.01"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("E38").Value = "  -4.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0150"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -2.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.782"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -9.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.495"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -4.27%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.755"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -2.49%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "93.21"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -5.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.06"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -3.00%  "
$ws.Range("D45").Value = "1.705.57"
$ws.Range("E45").Value = "  -3.60%  "
$ws.Range("D46").Value = "0.0₆0107"
$ws.Range("E46").Value = "  -7.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "52.90"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -3.54%  "
$ws.Range("E48").Value = "  -3.50%  "
$ws.Range("E49").Value = "  -4.74%  "
$ws.Range("E50").Value = "  -1.73%  "
$ws.Range("E51").Value = "  -0.30%  "
